$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961716033869"
$wb.Worksheets.Item(2).Name = "NB_TO-1650996173875386"
$wb.Worksheets.Item(3).Name = "RS_TO-1650996173875386"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961739394188"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961740114195"

# Sheet1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961715633795.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996171587417.csv"
$ws1.Range("B4").Value = "go_stims-1650996171587417.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961716033869.csv"

# Sheet2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509961734033883.csv"
$ws2.Range("B3").Value = "ZB-match_9-16509961720434115.csv"
$ws2.Range("B4").Value = "TB-16509961738514166.csv"
$ws2.Range("B5").Value = "OB-1650996172067381.csv"
$ws2.Range("B6").Value = "TB-16509961736914167.csv"
$ws2.Range("B7").Value = "TB-16509961736194248.csv"
$ws2.Range("B8").Value = "ZB-match_2-16509961716993792.csv"
$ws2.Range("B9").Value = "ZB-match_0-1650996171931389.csv"
$ws2.Range("B10").Value = "OB-1650996172603413.csv"

# Sheet4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961739074168.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961738833904.csv"
$ws4.Range("B4").Value = "MM_stims-16509961739234197.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961739074168.csv"
$ws4.Range("B6").Value = "MM_stims-16509961739394188.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961739234197.csv"

# Sheet5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509961739713817.csv"
$ws5.Range("B3").Value = "SAT_stims-1650996173955414.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961739394188.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650996173995436.csv"
